$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the "Result" (column C) values from "YES" to "NO" for all rows
# except Profile, which becomes "Yes" (lowercase variant).
$ws.Range("C2").Value = "NO"   # Registration
$ws.Range("C4").Value = "NO"   # Forgot Password
$ws.Range("C5").Value = "NO"   # Complete Course
$ws.Range("C6").Value = "NO"   # MC Distribution
$ws.Range("C7").Value = "NO"   # AccountInfo
$ws.Range("C8").Value = "Yes"  # Profile

# Update the selected cell in the sheet view
$ws.Range("A13").Select()
